# Update odds values in row 2
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.57
$ws.Range("I2").Value = 5.25
$ws.Range("L2").Value = 5.5
$ws.Range("O2").Value = 1.25
$ws.Range("P2").Value = 3.75
$ws.Range("Q2").Value = 1.8
$ws.Range("R2").Value = 2
$ws.Range("S2").Value = 1.33
$ws.Range("T2").Value = 3.25
$ws.Range("W2").Value = 7.5
$ws.Range("Z2").Value = 12
$ws.Range("AB2").Value = 23
$ws.Range("AC2").Value = 13
$ws.Range("AD2").Value = 8
$ws.Range("AE2").Value = 17
$ws.Range("AG2").Value = 251
$ws.Range("AN2").Value = 3.6
$ws.Range("AO2").Value = 8
$ws.Range("AT2").Value = 3.25
$ws.Range("AW2").Value = 7
$ws.Range("BB2").Value = 201

# Remove the last data row (row 4) entirely, which shifts the dimension
# of the sheet from A1:BD4 to A1:BD3
$ws.Rows.Item(4).Delete()
